{"js": "// Add a new row to the (last) table in the document body:\n//   column 1: \"16/01/2021\"\n//   column 2: \"mostrar comentarios de un monumento\" (same bullet-list\n//             formatting \u2014 style \"Prrafodelista\", numId 1 \u2014 as the rest\n//             of the \"Objetivo realizado\" entries).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The document has a single table (\"Fecha\" / \"Objetivo realizado\" log);\n// grab the last one defensively in case that ever changes.\nconst table = tables.items[tables.items.length - 1];\n\n// addRows inherits the paragraph/list formatting of the row above it for\n// each column, which is exactly how the existing rows in this table are\n// built (pStyle \"Prrafodelista\" + numPr ilvl 0 / numId 1 in column 2).\ntable.addRows(\"End\", 1, [\n  [\"16/01/2021\", \"mostrar comentarios de un monumento\"]\n]);\nawait context.sync();\n", "ps1": "# Add a new row to the (last) table in the document body:\n#   column 1: \"16/01/2021\"\n#   column 2: \"mostrar comentarios de un monumento\" (same bullet-list\n#             formatting - style \"Prrafodelista\", numId 1 - as the rest\n#             of the \"Objetivo realizado\" entries).\n$d = $word.ActiveDocument\n\n# The document has a single table (\"Fecha\" / \"Objetivo realizado\" log);\n# grab the last one defensively in case that ever changes.\n$table = $d.Tables.Item($d.Tables.Count)\n\n# Rows.Add() appends a new row at the end, inheriting each column's\n# paragraph/list formatting from the row above it - exactly how the\n# existing rows in this table are built (pStyle \"Prrafodelista\" +\n# numPr ilvl 0 / numId 1 in column 2).\n$newRow = $table.Rows.Add()\n$newRow.Cells(1).Range.Text = \"16/01/2021\"\n$newRow.Cells(2).Range.Text = \"mostrar comentarios de un monumento\"\n"}
